$d = $word.ActiveDocument

$replacements = @(
    @("705÷7=", "161÷8="),
    @("814÷5=", "966÷8="),
    @("978÷5=", "520÷8="),
    @("834÷4=", "651÷8="),
    @("851÷8=", "115÷2="),
    @("126÷2=", "723÷7="),
    @("906÷4=", "576÷6="),
    @("918÷7=", "725÷6="),
    @("794÷4=", "112÷2="),
    @("630÷2=", "537÷7="),
    @("428÷8=", "609÷7="),
    @("173÷8=", "495÷5="),
    @("394÷4=", "458÷2="),
    @("135÷6=", "481÷9="),
    @("629÷9=", "419÷7="),
    @("101÷4=", "523÷5="),
    @("411÷3=", "840÷3="),
    @("843÷2=", "183÷7="),
    @("350÷2=", "488÷3="),
    @("845÷4=", "849÷7="),
    @("376÷7=", "864÷4="),
    @("880÷5=", "291÷2="),
    @("840÷2=", "808÷4="),
    @("658÷5=", "432÷6="),
    @("826÷6=", "466÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
